# Edit script for 0709.docx: font fix + title/author/email swap + body/summary rewrite
$d = $word.ActiveDocument

# 1) Normalize font across the whole document body (TimesNewToman -> Times New Roman)
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"
$fullRange.Font.NameAscii = "Times New Roman"

# 2) Title (paragraph 1)
$p1 = $d.Paragraphs(1).Range
$d.Range($p1.Start, $p1.End - 1).Text = 'Shaping Tomorrow''s Minds: The Profound Impact of History on Our Future'

# 3) Author (paragraph 2)
$p2 = $d.Paragraphs(2).Range
$d.Range($p2.Start, $p2.End - 1).Text = 'Eleanor Richards'

# 4) Email (paragraph 3): runs were "l" "." "jireckova@astroyale" "." "edu"
#    target is "erichards@brightschools" "." "edu" (drop the astroyale runs)
$p3 = $d.Paragraphs(3).Range
$p3start = $p3.Start
$lenL = 1
$lenDot1 = 1
$lenAstro = ("jireckova@astroyale").Length
$lenDot2 = 1
$lenEdu = ("edu").Length
$runLend = $p3start + $lenL
$runDot1end = $runLend + $lenDot1
$runAstroEnd = $runDot1end + $lenAstro
$runDot2End = $runAstroEnd + $lenDot2
$runEduEnd = $runDot2End + $lenEdu
# delete back-to-front so earlier offsets stay valid
$d.Range($runAstroEnd, $runDot2End).Text = ""
$d.Range($runDot1end, $runAstroEnd).Text = ""
$d.Range($p3start, $runLend).Text = "erichards@brightschools"

# 5) Body paragraph (paragraph 5)
$p5 = $d.Paragraphs(5).Range
$d.Range($p5.Start, $p5.End - 1).Text = 'In the vast tapestry of human existence, few disciplines hold such profound sway over our present and future as history. It serves as the beacon illuminating our path, shedding light on the past to illuminate the path towards a better tomorrow. From the birth of great civilizations to the ebb and flow of empires, the study of history offers a panoramic view of human achievement and folly, imparting timeless lessons that shape our understanding of the world.' + ([char]11) + '' + ([char]11) + 'It is within the annals of history that we find the blueprint of our collective identity, a kaleidoscope of diverse cultures and traditions that have come together to create the rich mosaic of humanity. The study of these myriad threads that weave together the fabric of our existence fosters tolerance, understanding, and respect for the myriad ways in which civilizations and individuals have navigated the ebb and flow of time. History thus becomes a potent force for unity, steering us towards a future where differences are embraced rather than feared.' + ([char]11) + '' + ([char]11) + 'Moreover, history provides an unparalleled lens through which we can examine the intricacies of human nature. As we delve into the lives of influential figures from across time, we gain insights into the motivations, passions, and fears that drive our actions. This introspective journey helps us better comprehend our own place in the grand scheme of things, instilling empathy and compassion for the human condition. Through this prism, history guides us towards becoming more thoughtful decision-makers, whose actions are informed by a deep appreciation for both the complexity of the past and the potential of the future.'

# 6) Summary heading (paragraph 6) stays "Summary" - no change needed

# 7) Summary body (paragraph 7)
$p7 = $d.Paragraphs(7).Range
$d.Range($p7.Start, $p7.End - 1).Text = 'In this essay, we have explored the profound impact of history on our present and future. We have delved into the role of history in shaping our understanding of the world, inculcating tolerance and respect for diverse cultures, and providing insights into the complexities of human nature. History, therefore, serves as a vital tool for creating more thoughtful and compassionate citizens who are equipped with the knowledge and skills necessary to forge a better future for all.'

# 8) Append a new empty paragraph at the very end of the document
$endRng = $d.Range($d.Content.End, $d.Content.End)
$endRng.InsertParagraphAfter()
